# Auto-generated edit script: updates crypto price/volume table cells
# Applies the diff between before.xlsx and the refreshed crypto data export
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.616.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.965.71"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.39"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -7.26%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.972.57"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.113"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.33%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -8.25%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.479.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.39%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.636.09"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.61"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.965.81"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -6.25%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.12"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "379.49"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -7.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.90"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.62%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.98"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.48%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.083.38"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.59%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0925"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -10.26%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.99%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.72"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.37"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.94"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.56%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.28%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -9.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.91"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.405.66"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -10.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.93"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.10"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -8.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.664"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0590"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.41%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.92"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -9.91%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.58"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.22%  "
